$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a plain numeric-looking string into a "Price" cell (column D)
# while keeping it stored as TEXT (matching the source data, which is
# pre-formatted text like "5.43", not a real number). Forcing the
# NumberFormat to Text ("@") before the assignment stops Excel's automatic
# number conversion; resetting the Style back to "Normal" afterwards drops
# the now-unneeded explicit number format from the cell again so the cell
# keeps its original (unstyled) appearance.
function Set-TextPrice($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "62.923.66"
$ws.Range("E2").Value = "  -0.96%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.268.42"
$ws.Range("E3").Value = "  +0.40%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.08%  "

# Row 5 - BNB
Set-TextPrice "D5" "598.84"
$ws.Range("E5").Value = "  +0.51%  "

# Row 6 - Solana
Set-TextPrice "D6" "138.12"
$ws.Range("E6").Value = "  -2.38%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.11%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "3.266.29"
$ws.Range("E8").Value = "  +0.51%  "

# Row 9 - XRP
Set-TextPrice "D9" "0.510"
$ws.Range("E9").Value = "  -1.78%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -0.76%  "

# Row 11 - Toncoin
Set-TextPrice "D11" "5.43"
$ws.Range("E11").Value = "  +0.59%  "

# Row 12 - Cardano
Set-TextPrice "D12" "0.460"
$ws.Range("E12").Value = "  -1.48%  "

# Row 13 - ShibaInu
Set-TextPrice "D13" "0.0000241"
$ws.Range("E13").Value = "  -2.74%  "

# Row 14 - Avalanche
Set-TextPrice "D14" "33.90"
$ws.Range("E14").Value = "  -1.48%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.815.46"
$ws.Range("E15").Value = "  +0.70%  "

# Row 16 - TRON
$ws.Range("E16").Value = "  +1.07%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "3.273.48"
$ws.Range("E17").Value = "  +0.66%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "63.070.18"
$ws.Range("E18").Value = "  -0.74%  "

# Row 19 - Polkadot
Set-TextPrice "D19" "6.72"
$ws.Range("E19").Value = "  -0.93%  "

# Row 20 - BitcoinCash
Set-TextPrice "D20" "471.01"
$ws.Range("E20").Value = "  -1.68%  "

# Row 21 - Chainlink
Set-TextPrice "D21" "13.78"
$ws.Range("E21").Value = "  -3.25%  "

# Row 22 - Polygon
Set-TextPrice "D22" "0.724"
$ws.Range("E22").Value = "  -1.44%  "

# Row 23 - Uniswap
Set-TextPrice "D23" "7.84"
$ws.Range("E23").Value = "  -1.83%  "

# Row 24 - InternetComputer(DFINITY)
Set-TextPrice "D24" "13.65"
$ws.Range("E24").Value = "  +2.75%  "

# Row 25 - Litecoin
Set-TextPrice "D25" "84.11"
$ws.Range("E25").Value = "  +0.41%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.04%  "

# Row 27 - PancakeSwap
Set-TextPrice "D27" "2.73"
$ws.Range("E27").Value = "  -0.58%  "

# Row 28 - FirstDigitalUSD
Set-TextPrice "D28" "1.00"
$ws.Range("E28").Value = "  +0.02%  "

# Row 29 - NEARProtocol
Set-TextPrice "D29" "7.05"
$ws.Range("E29").Value = "  -2.39%  "

# Row 30 - RenderToken
Set-TextPrice "D30" "7.94"
$ws.Range("E30").Value = "  -2.04%  "

# Row 31 - ImmutableX
$ws.Range("E31").Value = "  -1.54%  "

# Row 32 - EthereumClassic
Set-TextPrice "D32" "28.06"
$ws.Range("E32").Value = "  +1.09%  "

# Row 33 - Hedera
Set-TextPrice "D33" "0.103"
$ws.Range("E33").Value = "  -3.53%  "

# Row 34 - Stacks
Set-TextPrice "D34" "2.46"
$ws.Range("E34").Value = "  -3.71%  "

# Row 35 - Mantle
$ws.Range("E35").Value = "  -1.22%  "

# Row 36 - Filecoin
Set-TextPrice "D36" "5.91"
$ws.Range("E36").Value = "  -0.62%  "

# Row 37 - OKB
Set-TextPrice "D37" "51.70"
$ws.Range("E37").Value = "  -2.05%  "

# Row 38 - PEPE
$ws.Range("D38").Value = "0.0₃0718"
$ws.Range("E38").Value = "  +0.37%  "

# Row 39 - VeChain
Set-TextPrice "D39" "0.0393"
$ws.Range("E39").Value = "  -0.32%  "

# Row 40 - Maker
$ws.Range("D40").Value = "3.081.72"
$ws.Range("E40").Value = "  +2.68%  "

# Row 41 - Bittensor
Set-TextPrice "D41" "420.63"
$ws.Range("E41").Value = "  -0.71%  "

# Row 42 - Kaspa
Set-TextPrice "D42" "0.116"
$ws.Range("E42").Value = "  +5.78%  "

# Row 43 - Cosmos
Set-TextPrice "D43" "8.19"
$ws.Range("E43").Value = "  -2.37%  "

# Row 44 - dogwifhat
Set-TextPrice "D44" "2.64"
$ws.Range("E44").Value = "  -4.70%  "

# Row 45 - TheGraph
Set-TextPrice "D45" "0.257"
$ws.Range("E45").Value = "  -3.09%  "

# Row 46 - Fetch.AI
Set-TextPrice "D46" "2.16"
$ws.Range("E46").Value = "  -1.01%  "

# Row 47 - USDe
$ws.Range("E47").Value = "  -0.07%  "

# Row 48 - was Arweave, now Monero
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextPrice "D48" "126.85"
$ws.Range("E48").Value = "  +3.29%  "

# Row 49 - was Monero, now Arweave
$ws.Range("B49").Value = "Arweave"
$ws.Range("C49").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextPrice "D49" "35.58"
$ws.Range("E49").Value = "  +5.93%  "

# Row 50 - InjectiveProtocol
Set-TextPrice "D50" "25.76"
$ws.Range("E50").Value = "  -0.72%  "

# Row 51 - Stellar
$ws.Range("E51").Value = "  -1.82%  "
